# 自动更新价格数据 2026-02-09 03:36:57
# Prepend a new day's row (2026-02-09) above the existing data, pushing all
# prior rows down by one. The new row carries the same price figures as
# every other row in this table (783.5 / 1112 / 3610).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows (old row 2.. ) down by inserting a blank row
# right below the header row.
$ws.Rows.Item(2).Insert()

# Force column A to stay plain text (matches the rest of the date column,
# which is stored as literal strings, not real dates) before writing the
# new date so Excel doesn't auto-convert "2026-02-09" into a date serial.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-09"
$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
